# Auto-generated: applies scheduled-runner market price refresh to Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 657.2727
$ws.Range("I2").Value = 403.75
$ws.Range("J2").Value = 1333.3334
$ws.Range("K2").Value = 403.75
$ws.Range("L2").Value = 1333.3334
$ws.Range("M2").Value = -290.75
$ws.Range("N2").Value = -1559.3334
$ws.Range("H5").Value = 170.66667
$ws.Range("I5").Value = 170.66667
$ws.Range("K5").Value = 170.66667
$ws.Range("M5").Value = -55.66667000000001
$ws.Range("H12").Value = 499
$ws.Range("I12").Value = 298.33334
$ws.Range("J12").Value = 699.6667
$ws.Range("K12").Value = 298.33334
$ws.Range("L12").Value = 699.6667
$ws.Range("M12").Value = -128.33334
$ws.Range("N12").Value = -1039.6667
$ws.Range("H18").Value = 1109
$ws.Range("I18").Value = 1109
$ws.Range("K18").Value = 1109
$ws.Range("M18").Value = -825
$ws.Range("H57").Value = 76332.664
$ws.Range("I57").Value = 69999
$ws.Range("J57").Value = 79499.5
$ws.Range("K57").Value = 209997
$ws.Range("L57").Value = 238498.5
$ws.Range("M57").Value = -209498
$ws.Range("N57").Value = -239496.5
$ws.Range("H74").Value = 3983.8572
$ws.Range("I74").Value = 3721.75
$ws.Range("K74").Value = 3721.75
$ws.Range("M74").Value = -2785.75
$ws.Range("H77").Value = 3983.8572
$ws.Range("I77").Value = 3721.75
$ws.Range("K77").Value = 18608.75
$ws.Range("M77").Value = -13928.75
$ws.Range("H86").Value = 28646270
$ws.Range("J86").Value = 66834800
$ws.Range("L86").Value = 66834800
$ws.Range("N86").Value = -66837046
$ws.Range("H89").Value = 28646270
$ws.Range("J89").Value = 66834800
$ws.Range("L89").Value = 334174000
$ws.Range("N89").Value = -334185232
$ws.Range("H98").Value = 11168.333
$ws.Range("I98").Value = 15752.5
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 15752.5
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -14254.5
$ws.Range("N98").Value = -4996
$ws.Range("H112").Value = 57071.61
$ws.Range("I112").Value = 112287
$ws.Range("K112").Value = 336861
$ws.Range("M112").Value = -335753
$ws.Range("H113").Value = 111115490
$ws.Range("I113").Value = 200003520
$ws.Range("J113").Value = 5450.25
$ws.Range("K113").Value = 200003520
$ws.Range("L113").Value = 5450.25
$ws.Range("M113").Value = -200000266
$ws.Range("N113").Value = -11958.25
$ws.Range("H122").Value = 11168.333
$ws.Range("I122").Value = 15752.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 47257.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -44807.5
$ws.Range("N122").Value = -10900
$ws.Range("H137").Value = 1842.2858
$ws.Range("I137").Value = 1300.375
$ws.Range("J137").Value = 3576.4
$ws.Range("K137").Value = 3901.125
$ws.Range("L137").Value = 10729.2
$ws.Range("M137").Value = -1351.125
$ws.Range("N137").Value = -15829.2
$ws.Range("H138").Value = 2312.0967
$ws.Range("J138").Value = 2793.7646
$ws.Range("L138").Value = 8381.293799999999
$ws.Range("N138").Value = -18661.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1151.05
$ws.Range("J2").Value = 997.8333
$ws.Range("L2").Value = 997.8333
$ws.Range("N2").Value = -1223.8333
$ws.Range("H32").Value = 4069.2163
$ws.Range("J32").Value = 3761
$ws.Range("L32").Value = 3761
$ws.Range("N32").Value = -4335
$ws.Range("H61").Value = 25644786
$ws.Range("I61").Value = 37039470
$ws.Range("K61").Value = 37039470
$ws.Range("M61").Value = -37039258
$ws.Range("H110").Value = 71500640
$ws.Range("I110").Value = 111167040
$ws.Range("J110").Value = 101129.8
$ws.Range("K110").Value = 111167040
$ws.Range("L110").Value = 101129.8
$ws.Range("M110").Value = -111164995
$ws.Range("N110").Value = -105219.8
$ws.Range("H116").Value = 1151.05
$ws.Range("J116").Value = 997.8333
$ws.Range("L116").Value = 997.8333
$ws.Range("N116").Value = -5585.8333
$ws.Range("H122").Value = 16669522
$ws.Range("I122").Value = 17546766
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 52640298
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -52637848
$ws.Range("N122").Value = -10600
$ws.Range("H132").Value = 71431140
$ws.Range("I132").Value = 83335820
$ws.Range("K132").Value = 250007460
$ws.Range("M132").Value = -250004930
$ws.Range("H136").Value = 25644786
$ws.Range("I136").Value = 37039470
$ws.Range("K136").Value = 111118410
$ws.Range("M136").Value = -111115860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1151.05
$ws.Range("J3").Value = 997.8333
$ws.Range("L3").Value = 997.8333
$ws.Range("N3").Value = -1225.8333
$ws.Range("H20").Value = 30498.25
$ws.Range("J20").Value = 7331.6665
$ws.Range("L20").Value = 7331.6665
$ws.Range("N20").Value = -7825.6665
$ws.Range("H86").Value = 3864.0625
$ws.Range("I86").Value = 5186.25
$ws.Range("J86").Value = 2541.875
$ws.Range("K86").Value = 5186.25
$ws.Range("L86").Value = 2541.875
$ws.Range("M86").Value = -4063.25
$ws.Range("N86").Value = -4787.875
$ws.Range("H89").Value = 3864.0625
$ws.Range("I89").Value = 5186.25
$ws.Range("J89").Value = 2541.875
$ws.Range("K89").Value = 25931.25
$ws.Range("L89").Value = 12709.375
$ws.Range("M89").Value = -20315.25
$ws.Range("N89").Value = -23941.375
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H134").Value = 6759127.5
$ws.Range("I134").Value = 8066701
$ws.Range("K134").Value = 24200103
$ws.Range("M134").Value = -24197568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2940.942
$ws.Range("I31").Value = 2156.7354
$ws.Range("K31").Value = 2156.7354
$ws.Range("M31").Value = -1861.7354
$ws.Range("H34").Value = 2940.942
$ws.Range("I34").Value = 2156.7354
$ws.Range("K34").Value = 2156.7354
$ws.Range("M34").Value = -1954.7354
$ws.Range("H58").Value = 2908
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2908
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 2908
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -3314
$ws.Range("H132").Value = 3209.0557
$ws.Range("I132").Value = 3126
$ws.Range("K132").Value = 9378
$ws.Range("M132").Value = -6848
$ws.Range("H136").Value = 2908
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2908
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 8724
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -13824
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 68064.336
$ws.Range("J141").Value = 68064.336
$ws.Range("L141").Value = 68064.336
$ws.Range("N141").Value = -78424.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1174.1875
$ws.Range("I8").Value = 1174.1875
$ws.Range("K8").Value = 3522.5625
$ws.Range("M8").Value = -3383.5625
$ws.Range("H131").Value = 8356.323
$ws.Range("I131").Value = 1110.4286
$ws.Range("J131").Value = 10234.889
$ws.Range("K131").Value = 3331.2858
$ws.Range("L131").Value = 30704.667
$ws.Range("M131").Value = 1708.7142
$ws.Range("N131").Value = -40784.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1877.7
$ws.Range("I80").Value = 2024.091
$ws.Range("J80").Value = 1698.7778
$ws.Range("K80").Value = 2024.091
$ws.Range("L80").Value = 1698.7778
$ws.Range("M80").Value = -1026.091
$ws.Range("N80").Value = -3694.7778
$ws.Range("H83").Value = 1877.7
$ws.Range("I83").Value = 2024.091
$ws.Range("J83").Value = 1698.7778
$ws.Range("K83").Value = 10120.455
$ws.Range("L83").Value = 8493.889000000001
$ws.Range("M83").Value = -5128.455
$ws.Range("N83").Value = -18477.889
$ws.Range("H92").Value = 21687.625
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 21687.625
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 21687.625
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -25431.625
$ws.Range("H97").Value = 1487.375
$ws.Range("I97").Value = 1668.625
$ws.Range("J97").Value = 1306.125
$ws.Range("K97").Value = 1668.625
$ws.Range("L97").Value = 1306.125
$ws.Range("M97").Value = -1172.625
$ws.Range("N97").Value = -2298.125
$ws.Range("H102").Value = 1902.8
$ws.Range("I102").Value = 1902.8
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1902.8
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -280.8
$ws.Range("N102").ClearContents()
$ws.Range("H107").Value = 1524
$ws.Range("I107").Value = 449.5
$ws.Range("J107").Value = 2598.5
$ws.Range("K107").Value = 449.5
$ws.Range("L107").Value = 2598.5
$ws.Range("M107").Value = 1470.5
$ws.Range("N107").Value = -6438.5
$ws.Range("H122").Value = 3086.6
$ws.Range("I122").Value = 2809.4285
$ws.Range("J122").Value = 3733.3333
$ws.Range("K122").Value = 8428.2855
$ws.Range("L122").Value = 11199.9999
$ws.Range("M122").Value = -5978.2855
$ws.Range("N122").Value = -16099.9999
$ws.Range("H123").Value = 75625
$ws.Range("J123").Value = 75625
$ws.Range("L123").Value = 75625
$ws.Range("N123").Value = -80525
$ws.Range("H132").Value = 3330.1
$ws.Range("I132").Value = 2700.7058
$ws.Range("J132").Value = 4153.154
$ws.Range("K132").Value = 8102.117400000001
$ws.Range("L132").Value = 12459.462
$ws.Range("M132").Value = -5572.117400000001
$ws.Range("N132").Value = -17519.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 83335704
$ws.Range("I7").Value = 100001850
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 100001850
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -100001738
$ws.Range("N7").Value = -5224
$ws.Range("H22").Value = 812.4
$ws.Range("I22").Value = 790.5
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 790.5
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -495.5
$ws.Range("N22").Value = -1490
$ws.Range("H27").Value = 812.4
$ws.Range("I27").Value = 790.5
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 790.5
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -683.5
$ws.Range("N27").Value = -1114
$ws.Range("H30").Value = 1696
$ws.Range("I30").Value = 650.5
$ws.Range("J30").Value = 2393
$ws.Range("K30").Value = 650.5
$ws.Range("L30").Value = 2393
$ws.Range("M30").Value = -542.5
$ws.Range("N30").Value = -2609
$ws.Range("H82").Value = 844.8182
$ws.Range("I82").Value = 783.3333
$ws.Range("J82").Value = 918.6
$ws.Range("K82").Value = 783.3333
$ws.Range("L82").Value = 918.6
$ws.Range("M82").Value = -422.3333
$ws.Range("N82").Value = -1640.6
$ws.Range("H85").Value = 844.8182
$ws.Range("I85").Value = 783.3333
$ws.Range("J85").Value = 918.6
$ws.Range("K85").Value = 783.3333
$ws.Range("L85").Value = 918.6
$ws.Range("M85").Value = 464.6667
$ws.Range("N85").Value = -3414.6
$ws.Range("H122").Value = 3843
$ws.Range("I122").Value = 2849.4443
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 8548.332900000001
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -6098.332900000001
$ws.Range("N122").Value = -20900.0005
$ws.Range("H126").Value = 83335704
$ws.Range("I126").Value = 100001850
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 300005550
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -300003080
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 2617.0833
$ws.Range("I132").Value = 2488.2
$ws.Range("J132").Value = 3261.5
$ws.Range("K132").Value = 7464.599999999999
$ws.Range("L132").Value = 9784.5
$ws.Range("M132").Value = -4934.599999999999
$ws.Range("N132").Value = -14844.5
$ws.Range("H136").Value = 2556.6667
$ws.Range("I136").Value = 2037.9412
$ws.Range("J136").Value = 3438.5
$ws.Range("K136").Value = 6113.8236
$ws.Range("L136").Value = 10315.5
$ws.Range("M136").Value = -3563.8236
$ws.Range("N136").Value = -15415.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1705532.1
$ws.Range("J62").Value = 5623.75
$ws.Range("L62").Value = 5623.75
$ws.Range("N62").Value = -6871.75
$ws.Range("H65").Value = 1705532.1
$ws.Range("J65").Value = 5623.75
$ws.Range("L65").Value = 28118.75
$ws.Range("N65").Value = -34358.75
$ws.Range("H81").Value = 8338832
$ws.Range("I81").Value = 2812.875
$ws.Range("J81").Value = 25010870
$ws.Range("K81").Value = 5625.75
$ws.Range("L81").Value = 50021740
$ws.Range("M81").Value = -4564.75
$ws.Range("N81").Value = -50023862
$ws.Range("H84").Value = 8338832
$ws.Range("I84").Value = 2812.875
$ws.Range("J84").Value = 25010870
$ws.Range("K84").Value = 28128.75
$ws.Range("L84").Value = 250108700
$ws.Range("M84").Value = -22824.75
$ws.Range("N84").Value = -250119308
$ws.Range("H126").Value = 1683
$ws.Range("I126").Value = 1519.6
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 4558.799999999999
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -2088.799999999999
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 1620.1555
$ws.Range("I132").Value = 1428.1621
$ws.Range("J132").Value = 2508.125
$ws.Range("K132").Value = 4284.4863
$ws.Range("L132").Value = 7524.375
$ws.Range("M132").Value = -1754.4863
$ws.Range("N132").Value = -12584.375
